$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Binding Value Set for IVL_TS.operator (row 5, column Z)
$elements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Adjust column width to reflect the wider content (bestFit recalculation)
$elements.Columns.Item(26).ColumnWidth = 50.3
